$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INDIA")
$ws.Range("Z1").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("Z1").EntireColumn.Delete()
